$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = -7.982800000000003
$ws.Range("A9").Value = -22.01089999999999
$ws.Range("D9").Value = -7.833699999999996
$ws.Range("D11").Value = -7.690699999999996
$ws.Range("A13").Value = -22.23249999999999
$ws.Range("A16").Value = -22.0493
$ws.Range("D16").Value = -8.483799999999995
$ws.Range("A18").Value = -22.26100000000001
$ws.Range("A20").Value = -21.43779999999998
$ws.Range("D23").Value = -8.1652
$ws.Range("D24").Value = -8.106199999999994
$ws.Range("A26").Value = -21.51449999999998
$ws.Range("D26").Value = -7.286299999999995
$ws.Range("A27").Value = -21.99099999999999
$ws.Range("A29").Value = -21.58899999999997
$ws.Range("D34").Value = -7.988900000000004
$ws.Range("A35").Value = -19.4777
$ws.Range("D35").Value = -7.759400000000003
$ws.Range("A36").Value = -20.80369999999998
$ws.Range("D44").Value = -7.869599999999997
$ws.Range("A45").Value = -21.71449999999998
$ws.Range("D48").Value = -7.795399999999997
$ws.Range("D49").Value = -8.059600000000001
$ws.Range("D52").Value = -7.756999999999998
$ws.Range("A55").Value = -22.22729999999999
$ws.Range("A57").Value = -22.39700000000001
$ws.Range("D66").Value = -7.273099999999999
$ws.Range("D67").Value = -7.4566
$ws.Range("A69").Value = -21.59119999999998
$ws.Range("D73").Value = -8.797699999999999
$ws.Range("A76").Value = -19.69089999999998
$ws.Range("A78").Value = -19.55419999999999
$ws.Range("D78").Value = -7.592800000000003
$ws.Range("D80").Value = -8.531599999999999
$ws.Range("A82").Value = -21.73719999999999
$ws.Range("A83").Value = -21.68539999999999
$ws.Range("D91").Value = -8.1351
$ws.Range("A93").Value = -21.31050000000003
$ws.Range("A97").Value = -21.52789999999998
$ws.Range("D97").Value = -7.378899999999996
$ws.Range("D99").Value = -8.099700000000002
$ws.Range("D104").Value = -7.789300000000001
